$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7262902855873108
$ws.Range("B1").Value = 1.384388208389282
$ws.Range("C1").Value = 4.390127182006836
$ws.Range("D1").Value = 1.824201464653015
$ws.Range("E1").Value = 1.045427560806274
